# hainan_covid19.xlsx - "data updated on Aug.12"
# Append three new days of Hainan COVID-19 data (Aug 9-11, 2022) to Sheet1,
# then move the view/selection the way the author left it (scrolled down,
# active cell D20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows: date (serial), hncon, hnasy, hnasytocon
$newRows = @(
    @(44782, 285, 285, 14),
    @(44783, 559, 805, 9),
    @(44784, 595, 614, 1)
)

$startRow = 10
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $values = $newRows[$i]

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Value = $values[0]
    $dateCell.NumberFormat = "M/D/YYYY"

    $ws.Cells.Item($r, 2).Value = $values[1]
    $ws.Cells.Item($r, 3).Value = $values[2]
    $ws.Cells.Item($r, 4).Value = $values[3]
}

# Scroll the view down to the newly added rows and leave the active
# selection where the author left it.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D20").Select()
